$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the data field width from 10 bits to 16 bits
$ws.Range("E1").Value = "data 16 bits"

# Move the selection to match the author's final cursor position
$ws.Range("G3:G4").Select()
